$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.792.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.074.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.070.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.80%  "

$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.589.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.866.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.085.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("E23").Value = "  +3.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("E34").Value = "  -3.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "447.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0814"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.964.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("E43").Value = "  -4.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("E51").Value = "  -0.62%  "
